# Insert a new row before row 105 (shifts rows 105:228 down to 106:229),
# then populate the new row 105 with a fresh data record (same
# Mercado/Region/Categoria/Variedad/Calidad/Unidad/Origen as the record
# that used to sit at row 105, but new Fecha/Volumen/Precio values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(105).Insert()

$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 44705
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = 100112043
$ws.Range("G105").Value = "Pepino dulce"
$ws.Range("H105").Value = "Cultivar IV Región"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 400
$ws.Range("K105").Value = 19000
$ws.Range("L105").Value = 20000
$ws.Range("M105").Value = 19500
$ws.Range("N105").Value = "$/bandeja 18 kilos"
$ws.Range("O105").Value = "Provincia de Limarí"
$ws.Range("P105").Value = 1083
$ws.Range("Q105").Value = 18
$ws.Range("R105").Value = "Hortaliza"
